$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 223, shifting existing rows 223:254 down to 224:255
$ws.Rows.Item(223).Insert()

# Fill in the new row 223 with its data
$ws.Cells.Item(223, 1).Value = 3
$ws.Cells.Item(223, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(223, 3).Value = "Coquimbo"
$ws.Cells.Item(223, 4).Value = 44491
$ws.Cells.Item(223, 5).Value = 5
$ws.Cells.Item(223, 6).Value = 100112017
$ws.Cells.Item(223, 7).Value = "Apio"
$ws.Cells.Item(223, 8).Value = "Americana (o)"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 180
$ws.Cells.Item(223, 11).Value = 9000
$ws.Cells.Item(223, 12).Value = 9000
$ws.Cells.Item(223, 13).Value = 9000
$ws.Cells.Item(223, 14).Value = "`$/docena de matas"
$ws.Cells.Item(223, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(223, 16).Value = 1500
$ws.Cells.Item(223, 17).Value = 6
$ws.Cells.Item(223, 18).Value = "Hortaliza"
